# Reorders the "Recorded By" (column G) values so that a leading "System, "
# prefix is moved to the end of the comma-separated list instead, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com" -> "backup@backdoor.com, System"
#   "System, system, backup@backdoor.com" -> "system, backup@backdoor.com, System"
# Rows whose "Recorded By" value also contains "admin@admin.com" are left
# untouched, as are rows that don't start with the "System, " prefix.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$prefix = "System, "

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value

    if ($text.StartsWith($prefix) -and -not $text.Contains("admin@admin.com")) {
        $rest = $text.Substring($prefix.Length)
        $cell.Value = "$rest, System"
    }
}
